# Order Creation: Show order and Add On Order
# Update environment target from test18 -> test14, update the
# active cell selection, and keep the A2 hyperlink's stale display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the old (test18) display text on the A2 hyperlink before the
# underlying cell text is changed to test14, so the saved hyperlink keeps
# display="https://test18.cliotest.com/backoffice/control/main" while its
# relationship target is unchanged.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "https://test18.cliotest.com/backoffice/control/main"
    }
}

# Update the environment-specific cell values from test18 to test14.
$ws.Range("A2").Value = "https://test14.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test14.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test14.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest14"
$ws.Range("G2").Value = "test14"
$ws.Range("K2").Value = "test14"

# Move the active selection to C13 (also clears the old E1 scroll position).
$ws.Range("C13").Select()
